$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21: Book and a Hard Place / Engraved Hard Leather Grimoire
$ws.Cells.Item(21, 8).Value = 20266
$ws.Cells.Item(21, 9).Value = 800
$ws.Cells.Item(21, 11).Value = 800
$ws.Cells.Item(21, 13).Value = -332

# Row 23: There's Something about Bury / Hard Leather Grimoire
$ws.Cells.Item(23, 8).Value = 20266
$ws.Cells.Item(23, 9).Value = 800
$ws.Cells.Item(23, 11).Value = 800
$ws.Cells.Item(23, 13).Value = -566

# Row 53: No Accounting for Waste / Enchanted Electrum Ink
$ws.Cells.Item(53, 8).Value = 213.15384
$ws.Cells.Item(53, 9).Value = 119
$ws.Cells.Item(53, 10).Value = 425
$ws.Cells.Item(53, 11).Value = 119
$ws.Cells.Item(53, 12).Value = 425
$ws.Cells.Item(53, 13).Value = 518
$ws.Cells.Item(53, 14).Value = -1699

# Row 69: Steeling the Knife, Steeling the Mind / Grade 1 Mind Dissolvent
$ws.Cells.Item(69, 8).Value = 3812.8572
$ws.Cells.Item(69, 9).Value = 4188
$ws.Cells.Item(69, 10).Value = 2875
$ws.Cells.Item(69, 11).Value = 12564
$ws.Cells.Item(69, 12).Value = 8625
$ws.Cells.Item(69, 13).Value = -11690
$ws.Cells.Item(69, 14).Value = -10373

# Row 72: Surgical Substitution (L) / Grade 1 Mind Dissolvent
$ws.Cells.Item(72, 8).Value = 3812.8572
$ws.Cells.Item(72, 9).Value = 4188
$ws.Cells.Item(72, 10).Value = 2875
$ws.Cells.Item(72, 11).Value = 37692
$ws.Cells.Item(72, 12).Value = 25875
$ws.Cells.Item(72, 13).Value = -33324
$ws.Cells.Item(72, 14).Value = -34611

# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Cells.Item(112, 8).Value = 3133.6155
$ws.Cells.Item(112, 9).Value = 2050
$ws.Cells.Item(112, 10).Value = 3192.1892
$ws.Cells.Item(112, 11).Value = 6150
$ws.Cells.Item(112, 12).Value = 9576.567599999998
$ws.Cells.Item(112, 13).Value = -5042
$ws.Cells.Item(112, 14).Value = -11792.5676

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Cells.Item(132, 8).Value = 3911.25
$ws.Cells.Item(132, 9).Value = 3620.7896
$ws.Cells.Item(132, 11).Value = 10862.3688
$ws.Cells.Item(132, 13).Value = -8332.3688

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Cells.Item(138, 8).Value = 5428.6177
$ws.Cells.Item(138, 9).Value = 9724.25
$ws.Cells.Item(138, 10).Value = 4855.8667
$ws.Cells.Item(138, 11).Value = 29172.75
$ws.Cells.Item(138, 12).Value = 14567.6001
$ws.Cells.Item(138, 13).Value = -24032.75
$ws.Cells.Item(138, 14).Value = -24847.6001

$ws = $wb.Worksheets.Item("ARM")
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Cells.Item(74, 8).Value = 985.5238000000001
$ws.Cells.Item(74, 9).Value = 833.069
$ws.Cells.Item(74, 11).Value = 833.069
$ws.Cells.Item(74, 13).Value = 40.93100000000004

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Cells.Item(77, 8).Value = 985.5238000000001
$ws.Cells.Item(77, 9).Value = 833.069
$ws.Cells.Item(77, 11).Value = 4165.344999999999
$ws.Cells.Item(77, 13).Value = 202.6550000000007

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Cells.Item(132, 8).Value = 4459.528
$ws.Cells.Item(132, 9).Value = 3557.5925
$ws.Cells.Item(132, 10).Value = 7165.3335
$ws.Cells.Item(132, 11).Value = 10672.7775
$ws.Cells.Item(132, 12).Value = 21496.0005
$ws.Cells.Item(132, 13).Value = -8142.7775
$ws.Cells.Item(132, 14).Value = -26556.0005

# Row 134: Brace for More Vambraces / Ruthenium Vambraces of Maiming
$ws.Cells.Item(134, 8).Value = 74419.336
$ws.Cells.Item(134, 10).Value = 74419.336
$ws.Cells.Item(134, 12).Value = 74419.336
$ws.Cells.Item(134, 14).Value = -84559.336

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Cells.Item(31, 8).Value = 1211.8334
$ws.Cells.Item(31, 9).Value = 1380
$ws.Cells.Item(31, 11).Value = 1380
$ws.Cells.Item(31, 13).Value = -1085

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Cells.Item(34, 8).Value = 1211.8334
$ws.Cells.Item(34, 9).Value = 1380
$ws.Cells.Item(34, 11).Value = 1380
$ws.Cells.Item(34, 13).Value = -1178

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Cells.Item(58, 8).Value = 1117.1077
$ws.Cells.Item(58, 9).Value = 946.7954999999999
$ws.Cells.Item(58, 10).Value = 1473.9524
$ws.Cells.Item(58, 11).Value = 946.7954999999999
$ws.Cells.Item(58, 12).Value = 1473.9524
$ws.Cells.Item(58, 13).Value = -743.7954999999999
$ws.Cells.Item(58, 14).Value = -1879.9524

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Cells.Item(122, 8).Value = 1425
$ws.Cells.Item(122, 9).Value = 1376.1875
$ws.Cells.Item(122, 10).Value = 1496
$ws.Cells.Item(122, 11).Value = 4128.5625
$ws.Cells.Item(122, 12).Value = 4488
$ws.Cells.Item(122, 13).Value = -1678.5625
$ws.Cells.Item(122, 14).Value = -9388

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Cells.Item(136, 8).Value = 1117.1077
$ws.Cells.Item(136, 9).Value = 946.7954999999999
$ws.Cells.Item(136, 10).Value = 1473.9524
$ws.Cells.Item(136, 11).Value = 2840.3865
$ws.Cells.Item(136, 12).Value = 4421.857199999999
$ws.Cells.Item(136, 13).Value = -290.3864999999996
$ws.Cells.Item(136, 14).Value = -9521.857199999999

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Cells.Item(5, 8).Value = 1056.6923
$ws.Cells.Item(5, 9).Value = 343.73914
$ws.Cells.Item(5, 10).Value = 2081.5625
$ws.Cells.Item(5, 11).Value = 1031.21742
$ws.Cells.Item(5, 12).Value = 6244.6875
$ws.Cells.Item(5, 13).Value = -919.2174199999999
$ws.Cells.Item(5, 14).Value = -6468.6875

# Row 7: It's Always Sunny in Vylbrand / Raisins
$ws.Cells.Item(7, 8).Value = 665.4545000000001
$ws.Cells.Item(7, 9).Value = 386.66666
$ws.Cells.Item(7, 10).Value = 1000
$ws.Cells.Item(7, 11).Value = 1159.99998
$ws.Cells.Item(7, 12).Value = 3000
$ws.Cells.Item(7, 13).Value = -1047.99998
$ws.Cells.Item(7, 14).Value = -3224

# Row 34: Fever Pitch / Chamomile Tea
$ws.Cells.Item(34, 8).Value = 20000506
$ws.Cells.Item(34, 10).Value = 27778386
$ws.Cells.Item(34, 12).Value = 83335158
$ws.Cells.Item(34, 14).Value = -83335326

# Row 39: Bloody Good Tart, This / Blood Currant Tart
$ws.Cells.Item(39, 8).Value = 4612.5
$ws.Cells.Item(39, 10).Value = 4612.5
$ws.Cells.Item(39, 12).Value = 13837.5
$ws.Cells.Item(39, 14).Value = -14425.5

# Row 55: Pagan Pastries / Pastry Fish
$ws.Cells.Item(55, 8).Value = 2016.6666
$ws.Cells.Item(55, 10).Value = 2016.6666
$ws.Cells.Item(55, 12).Value = 6049.9998
$ws.Cells.Item(55, 14).Value = -6403.9998

# Row 68: Such a Butter Face / Fermented Butter
$ws.Cells.Item(68, 8).Value = 902.8461
$ws.Cells.Item(68, 9).Value = 749
$ws.Cells.Item(68, 10).Value = 949
$ws.Cells.Item(68, 11).Value = 2247
$ws.Cells.Item(68, 12).Value = 2847
$ws.Cells.Item(68, 13).Value = -1436
$ws.Cells.Item(68, 14).Value = -4469

# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Cells.Item(71, 8).Value = 902.8461
$ws.Cells.Item(71, 9).Value = 749
$ws.Cells.Item(71, 10).Value = 949
$ws.Cells.Item(71, 11).Value = 6741
$ws.Cells.Item(71, 12).Value = 8541
$ws.Cells.Item(71, 13).Value = -2685
$ws.Cells.Item(71, 14).Value = -16653

# Row 75: Breakfast of Champions / Emerald Soup
$ws.Cells.Item(75, 8).Value = 710.75
$ws.Cells.Item(75, 9).Value = 713
$ws.Cells.Item(75, 10).Value = 710
$ws.Cells.Item(75, 11).Value = 2139
$ws.Cells.Item(75, 12).Value = 2130
$ws.Cells.Item(75, 13).Value = -1141
$ws.Cells.Item(75, 14).Value = -4126

# Row 78: Emerald Soup for the Soul (L) / Emerald Soup
$ws.Cells.Item(78, 8).Value = 710.75
$ws.Cells.Item(78, 9).Value = 713
$ws.Cells.Item(78, 10).Value = 710
$ws.Cells.Item(78, 11).Value = 6417
$ws.Cells.Item(78, 12).Value = 6390
$ws.Cells.Item(78, 13).Value = -1425
$ws.Cells.Item(78, 14).Value = -16374

# Row 92: Oh No Udon / Gyr Abanian Flour
$ws.Cells.Item(92, 8).Value = 646.8946999999999
$ws.Cells.Item(92, 9).Value = 595.2
$ws.Cells.Item(92, 10).Value = 704.3333
$ws.Cells.Item(92, 11).Value = 1785.6
$ws.Cells.Item(92, 12).Value = 2112.9999
$ws.Cells.Item(92, 13).Value = -537.6000000000001
$ws.Cells.Item(92, 14).Value = -4608.9999

# Row 122: Salt of the North / Northern Sea Salt
$ws.Cells.Item(122, 8).Value = 7696.143
$ws.Cells.Item(122, 9).Value = 324.8
$ws.Cells.Item(122, 10).Value = 26124.5
$ws.Cells.Item(122, 11).Value = 2923.2
$ws.Cells.Item(122, 12).Value = 235120.5
$ws.Cells.Item(122, 13).Value = -473.2000000000003
$ws.Cells.Item(122, 14).Value = -240020.5

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Cells.Item(135, 8).Value = 1056.6923
$ws.Cells.Item(135, 9).Value = 343.73914
$ws.Cells.Item(135, 10).Value = 2081.5625
$ws.Cells.Item(135, 11).Value = 3093.65226
$ws.Cells.Item(135, 12).Value = 18734.0625
$ws.Cells.Item(135, 13).Value = -558.6522600000003
$ws.Cells.Item(135, 14).Value = -23804.0625

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Cells.Item(80, 8).Value = 787491.7
$ws.Cells.Item(80, 9).Value = 1004766.7
$ws.Cells.Item(80, 10).Value = 135666.67
$ws.Cells.Item(80, 11).Value = 1004766.7
$ws.Cells.Item(80, 12).Value = 135666.67
$ws.Cells.Item(80, 13).Value = -1003768.7
$ws.Cells.Item(80, 14).Value = -137662.67

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Cells.Item(83, 8).Value = 787491.7
$ws.Cells.Item(83, 9).Value = 1004766.7
$ws.Cells.Item(83, 10).Value = 135666.67
$ws.Cells.Item(83, 11).Value = 5023833.5
$ws.Cells.Item(83, 12).Value = 678333.3500000001
$ws.Cells.Item(83, 13).Value = -5018841.5
$ws.Cells.Item(83, 14).Value = -688317.3500000001

# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Cells.Item(107, 8).Value = 797.6
$ws.Cells.Item(107, 9).Value = 603.25
$ws.Cells.Item(107, 10).Value = 1019.7143
$ws.Cells.Item(107, 11).Value = 603.25
$ws.Cells.Item(107, 12).Value = 1019.7143
$ws.Cells.Item(107, 13).Value = 1316.75
$ws.Cells.Item(107, 14).Value = -4859.7143

# Row 132: On Board for Lar / Lar Ingot
$ws.Cells.Item(132, 8).Value = 3892.2144
$ws.Cells.Item(132, 9).Value = 3721.4443
$ws.Cells.Item(132, 10).Value = 4199.6
$ws.Cells.Item(132, 11).Value = 11164.3329
$ws.Cells.Item(132, 12).Value = 12598.8
$ws.Cells.Item(132, 13).Value = -8634.332900000001
$ws.Cells.Item(132, 14).Value = -17658.8

$ws = $wb.Worksheets.Item("WVR")
# Row 93: What Guides Want / Bloodhempen Doublet of Crafting
$ws.Cells.Item(93, 8).Value = 53097.25
$ws.Cells.Item(93, 10).Value = 53097.25
$ws.Cells.Item(93, 12).Value = 53097.25
$ws.Cells.Item(93, 14).Value = -58089.25

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Cells.Item(132, 8).Value = 12684876
$ws.Cells.Item(132, 9).Value = 4439.8335
$ws.Cells.Item(132, 10).Value = 26518080
$ws.Cells.Item(132, 11).Value = 13319.5005
$ws.Cells.Item(132, 12).Value = 79554240
$ws.Cells.Item(132, 13).Value = -10789.5005
$ws.Cells.Item(132, 14).Value = -79559300

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Cells.Item(136, 8).Value = 2043.174
$ws.Cells.Item(136, 9).Value = 1745.5676
$ws.Cells.Item(136, 10).Value = 3266.6667
$ws.Cells.Item(136, 11).Value = 5236.7028
$ws.Cells.Item(136, 12).Value = 9800.000100000001
$ws.Cells.Item(136, 13).Value = -2686.7028
$ws.Cells.Item(136, 14).Value = -14900.0001
